$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 ("Save") - reuse the same formatting as the other
# header cells (bold/centered/bordered style) by copying from G1.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats

# New data column H2:H6, all zeros, unstyled like the other number cells.
$ws.Range("H2:H6").Value = 0
